# Weekly update: insert the latest "Albahaca" price-record at the top of
# the data table (row 10, just under the most recent existing entry),
# pushing all prior observations down by one row.
#
# Before: data rows were 2..59 (dimension A1:R59).
# After : a brand-new observation is inserted as row 10, all former rows
#         10..59 shift down to 11..60 (dimension A1:R60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10:59 down to 11:60, leaving a blank row 10 for the new entry.
$ws.Rows("10:10").Insert()

# Populate the new row 10 with the latest weekly observation.
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44425
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112052
$ws.Range("G10").Value = "Albahaca"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 90
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("N10").Value = "$/paquete"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 7000
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
